# Generate Report for Handback
# Row 7 (the f4d16b7e-8558-4edd-a02a-f84d975295e9 entry) finished its
# handback cycle, so the "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" / "Error Detail" columns (I, J, K, P) on both
# the "zh-cn" and "de-de" report sheets now get populated for that row.

$wb = $excel.ActiveWorkbook

$targetMdName = "f4d16b7e-8558-4edd-a02a-f84d975295e9.md"

$msg = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cbac89161ce2bcead9cce13cad90cc9d866134c0/e2e/f4d16b7e-8558-4edd-a02a-f84d975295e9.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/64a67f99ef374bb98259e5a329249294a6c14497/e2e/f4d16b7e-8558-4edd-a02a-f84d975295e9.md."

# ---------------------------------------------------------------------
# zh-cn sheet, row 7
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$iZh = $wsZh.Cells.Item(7, 9)
$iZh.Value = $targetMdName
$iZh.Style = "Hyperlink"
$wsZh.Hyperlinks.Add($iZh, "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/64a67f99ef374bb98259e5a329249294a6c14497/e2e/f4d16b7e-8558-4edd-a02a-f84d975295e9.md", $null, $null, $targetMdName) | Out-Null

$wsZh.Cells.Item(7, 10).Value = "f4d16b7e-8558-4edd-a02a-f84d975295e9.e84f954bc4977eca39021ffcdc0d5a44dcde57de.zh-cn.xlf"
$wsZh.Cells.Item(7, 11).Value = $msg
$wsZh.Cells.Item(7, 16).Value = "2016-08-27 19:01:46"

# ---------------------------------------------------------------------
# de-de sheet, row 7
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$iDe = $wsDe.Cells.Item(7, 9)
$iDe.Value = $targetMdName
$iDe.Style = "Hyperlink"
$wsDe.Hyperlinks.Add($iDe, "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/64a67f99ef374bb98259e5a329249294a6c14497/e2e/f4d16b7e-8558-4edd-a02a-f84d975295e9.md", $null, $null, $targetMdName) | Out-Null

$wsDe.Cells.Item(7, 10).Value = "f4d16b7e-8558-4edd-a02a-f84d975295e9.e84f954bc4977eca39021ffcdc0d5a44dcde57de.de-de.xlf"
$wsDe.Cells.Item(7, 11).Value = "2016-08-27 19:02:15"
$wsDe.Cells.Item(7, 16).Value = "2016-08-27 19:01:46"
